$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.012360031716525555
$ws.Range("C2").Value = 0.004973901901394129
$ws.Range("D2").Value = 0.003951451275497675
$ws.Range("E2").Value = 0.002330801449716091
$ws.Range("F2").Value = 0.0000031050494726514444
$ws.Range("J2").Value = 0.1261412650346756
$ws.Range("K2").Value = 1.4085742235183716
